# Apply "backup" column (R) + backfill 6 new monthly rows (266-271) to the
# AMBUJACEM.NS 1mo stock-history sheet, and zero out six stale
# `detect_structure` (Q) flags that were superseded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New header cell R1 = "backup" (copy Q1's bold/bordered header style)
# ---------------------------------------------------------------------
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("R1").Value = "backup"
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Column R ("backup") for existing data rows 2-265: mostly 0, except
#    rows 207 and 216 which back up their existing Q value of 2.
# ---------------------------------------------------------------------
$specialBackup = @{ 207 = 2; 216 = 2 }

for ($i = 2; $i -le 265; $i++) {
    if ($specialBackup.ContainsKey($i)) {
        $ws.Cells.Item($i, 18).Value = $specialBackup[$i]
    } else {
        $ws.Cells.Item($i, 18).Value = 0
    }
}

# ---------------------------------------------------------------------
# 3) Six rows whose `detect_structure` (Q) flag got cleared to 0 as part
#    of this pass (their old value is the one preserved above, if any).
# ---------------------------------------------------------------------
$qReset = 14, 29, 33, 36, 52, 53
foreach ($r in $qReset) {
    $ws.Cells.Item($r, 17).Value = 0
}

# ---------------------------------------------------------------------
# 4) Row 265's `isPivot` (O) flag flips from 0 to 2.
# ---------------------------------------------------------------------
$ws.Cells.Item(265, 15).Value = 2

# ---------------------------------------------------------------------
# 5) Append six freshly-pulled monthly candles (Jul-Dec 2024), rows
#    266-271. Columns: A Datetime, B Open, C High, D Low, E Close,
#    F Adj Close (blank/NaN for these newest rows), G Volume, H Year,
#    I Month, J Day, K Hour, L Minute, M Second, N Week, O isPivot,
#    P two_line_structure, Q detect_structure. R (backup) stays blank.
# ---------------------------------------------------------------------
$newRows = @(
    @(266, 45474, 668,               706.9500122070312, 656.2000122070312, 679.9500122070312, 75042770, 2024, 7,  1, 0, 0, 0, 27, 1, 0, 0),
    @(267, 45505, 680,               681,                600.7000122070312, 617.0499877929688, 70715449, 2024, 8,  1, 0, 0, 0, 31, 0, 0, 0),
    @(268, 45536, 620,               643.2999877929688, 598,                632.5499877929688, 41059526, 2024, 9,  1, 0, 0, 0, 35, 0, 0, 0),
    @(269, 45566, 630.2000122070312, 634.75,             545.2000122070312, 580.5499877929688, 84341651, 2024, 10, 1, 0, 0, 0, 40, 0, 0, 1),
    @(270, 45597, 582.9500122070312, 585.5,              453.0499877929688, 531.5,              93753069, 2024, 11, 1, 0, 0, 0, 44, 0, 0, 2),
    @(271, 45627, 532,               584.2000122070312, 530.0499877929688, 544.5999755859375,   39095578, 2024, 12, 1, 0, 0, 0, 48, 0, 0, 0)
)

foreach ($row in $newRows) {
    $r = $row[0]

    # Give the Datetime cell the same date style ("s=2") the rest of
    # column A uses, by copying the format down from the row above.
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Cells.Item($r, 1).Value  = $row[1]   # A Datetime
    $ws.Cells.Item($r, 2).Value  = $row[2]   # B Open
    $ws.Cells.Item($r, 3).Value  = $row[3]   # C High
    $ws.Cells.Item($r, 4).Value  = $row[4]   # D Low
    $ws.Cells.Item($r, 5).Value  = $row[5]   # E Close
    # F Adj Close intentionally left blank for these newest rows
    $ws.Cells.Item($r, 7).Value  = $row[6]   # G Volume
    $ws.Cells.Item($r, 8).Value  = $row[7]   # H Year
    $ws.Cells.Item($r, 9).Value  = $row[8]   # I Month
    $ws.Cells.Item($r, 10).Value = $row[9]   # J Day
    $ws.Cells.Item($r, 11).Value = $row[10]  # K Hour
    $ws.Cells.Item($r, 12).Value = $row[11]  # L Minute
    $ws.Cells.Item($r, 13).Value = $row[12]  # M Second
    $ws.Cells.Item($r, 14).Value = $row[13]  # N Week
    $ws.Cells.Item($r, 15).Value = $row[14]  # O isPivot
    $ws.Cells.Item($r, 16).Value = $row[15]  # P two_line_structure
    $ws.Cells.Item($r, 17).Value = $row[16]  # Q detect_structure
    # R backup intentionally left blank for these newest rows
}

Write-Host "Applied backup column + 6 appended rows"
